$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 172.35715
$ws.Range("I9").Value = 144.91667
$ws.Range("K9").Value = 144.91667
$ws.Range("M9").Value = 24.08332999999999
$ws.Range("H33").Value = 3519.125
$ws.Range("I33").Value = 3657.4783
$ws.Range("K33").Value = 3657.4783
$ws.Range("M33").Value = -3428.4783
$ws.Range("H40").Value = 4066.3809
$ws.Range("I40").Value = 1800.5714
$ws.Range("J40").Value = 5199.2856
$ws.Range("K40").Value = 1800.5714
$ws.Range("L40").Value = 5199.2856
$ws.Range("M40").Value = -1625.5714
$ws.Range("N40").Value = -5549.2856
$ws.Range("H74").Value = 7440.643
$ws.Range("I74").Value = 4174.5
$ws.Range("J74").Value = 7985
$ws.Range("K74").Value = 4174.5
$ws.Range("L74").Value = 7985
$ws.Range("M74").Value = -3238.5
$ws.Range("N74").Value = -9857
$ws.Range("H77").Value = 7440.643
$ws.Range("I77").Value = 4174.5
$ws.Range("J77").Value = 7985
$ws.Range("K77").Value = 20872.5
$ws.Range("L77").Value = 39925
$ws.Range("M77").Value = -16192.5
$ws.Range("N77").Value = -49285
$ws.Range("H80").Value = 621.1053000000001
$ws.Range("J80").Value = 676.6429000000001
$ws.Range("L80").Value = 2029.9287
$ws.Range("N80").Value = -4025.9287
$ws.Range("H83").Value = 621.1053000000001
$ws.Range("J83").Value = 676.6429000000001
$ws.Range("L83").Value = 6089.7861
$ws.Range("N83").Value = -16073.7861
$ws.Range("H86").Value = 1266
$ws.Range("J86").Value = 1499
$ws.Range("L86").Value = 1499
$ws.Range("N86").Value = -3745
$ws.Range("H89").Value = 1266
$ws.Range("J89").Value = 1499
$ws.Range("L89").Value = 7495
$ws.Range("N89").Value = -18727
$ws.Range("H92").Value = 1838.7
$ws.Range("I92").Value = 2111
$ws.Range("J92").Value = 749.5
$ws.Range("K92").Value = 2111
$ws.Range("L92").Value = 749.5
$ws.Range("M92").Value = -863
$ws.Range("N92").Value = -3245.5
$ws.Range("H106").Value = 2049.5
$ws.Range("I106").Value = 2159.4
$ws.Range("K106").Value = 2159.4
$ws.Range("M106").Value = -1528.4
$ws.Range("H112").Value = 22415.2
$ws.Range("J112").Value = 35365.668
$ws.Range("L112").Value = 106097.004
$ws.Range("N112").Value = -108313.004

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3537.5
$ws.Range("I61").Value = 2995.4285
$ws.Range("J61").Value = 4296.4
$ws.Range("K61").Value = 2995.4285
$ws.Range("L61").Value = 4296.4
$ws.Range("M61").Value = -2783.4285
$ws.Range("N61").Value = -4720.4
$ws.Range("H92").Value = 78961.164
$ws.Range("J92").Value = 78961.164
$ws.Range("L92").Value = 78961.164
$ws.Range("N92").Value = -83953.164
$ws.Range("H94").Value = 38666.332
$ws.Range("J94").Value = 38666.332
$ws.Range("L94").Value = 38666.332
$ws.Range("N94").Value = -40468.332
$ws.Range("H97").Value = 1618698
$ws.Range("I97").Value = 1618698
$ws.Range("K97").Value = 1618698
$ws.Range("M97").Value = -1618202
$ws.Range("H122").Value = 1304500.5
$ws.Range("I122").Value = 2481.5386
$ws.Range("J122").Value = 6946583
$ws.Range("K122").Value = 7444.6158
$ws.Range("L122").Value = 20839749
$ws.Range("M122").Value = -4994.6158
$ws.Range("N122").Value = -20844649
$ws.Range("H132").Value = 2435.9395
$ws.Range("I132").Value = 1613.84
$ws.Range("J132").Value = 5005
$ws.Range("K132").Value = 4841.52
$ws.Range("L132").Value = 15015
$ws.Range("M132").Value = -2311.52
$ws.Range("N132").Value = -20075
$ws.Range("H136").Value = 3537.5
$ws.Range("I136").Value = 2995.4285
$ws.Range("J136").Value = 4296.4
$ws.Range("K136").Value = 8986.2855
$ws.Range("L136").Value = 12889.2
$ws.Range("M136").Value = -6436.2855
$ws.Range("N136").Value = -17989.2

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4769985
$ws.Range("I86").Value = 6677451
$ws.Range("J86").Value = 1320.8334
$ws.Range("K86").Value = 6677451
$ws.Range("L86").Value = 1320.8334
$ws.Range("M86").Value = -6676328
$ws.Range("N86").Value = -3566.8334
$ws.Range("H89").Value = 4769985
$ws.Range("I89").Value = 6677451
$ws.Range("J89").Value = 1320.8334
$ws.Range("K89").Value = 33387255
$ws.Range("L89").Value = 6604.166999999999
$ws.Range("M89").Value = -33381639
$ws.Range("N89").Value = -17836.167
$ws.Range("H134").Value = 3801.4644
$ws.Range("I134").Value = 1611.2273
$ws.Range("J134").Value = 11832.333
$ws.Range("K134").Value = 4833.6819
$ws.Range("L134").Value = 35496.999
$ws.Range("M134").Value = -2298.6819
$ws.Range("N134").Value = -40566.999

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 33333
$ws.Range("J92").Value = 33333
$ws.Range("L92").Value = 33333
$ws.Range("N92").Value = -38325
$ws.Range("H93").Value = 31299.75
$ws.Range("I93").Value = 8483.166999999999
$ws.Range("K93").Value = 8483.166999999999
$ws.Range("M93").Value = -6611.166999999999
$ws.Range("H132").Value = 55976.113
$ws.Range("I132").Value = 36397.516
$ws.Range("K132").Value = 109192.548
$ws.Range("M132").Value = -106662.548
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 378.2963
$ws.Range("J2").Value = 455.73685
$ws.Range("L2").Value = 2734.4211
$ws.Range("N2").Value = -2960.4211
$ws.Range("H15").Value = 694.8333
$ws.Range("I15").Value = 692
$ws.Range("J15").Value = 696.25
$ws.Range("K15").Value = 2076
$ws.Range("L15").Value = 2088.75
$ws.Range("M15").Value = -1936
$ws.Range("N15").Value = -2368.75
$ws.Range("H86").Value = 444.6
$ws.Range("I86").Value = 399.5
$ws.Range("J86").Value = 474.66666
$ws.Range("K86").Value = 1198.5
$ws.Range("L86").Value = 1423.99998
$ws.Range("M86").Value = -12.5
$ws.Range("N86").Value = -3795.99998
$ws.Range("H89").Value = 444.6
$ws.Range("I89").Value = 399.5
$ws.Range("J89").Value = 474.66666
$ws.Range("K89").Value = 3595.5
$ws.Range("L89").Value = 4271.99994
$ws.Range("M89").Value = 2332.5
$ws.Range("N89").Value = -16127.99994
$ws.Range("H122").Value = 983.46155
$ws.Range("J122").Value = 989.75
$ws.Range("L122").Value = 8907.75
$ws.Range("N122").Value = -13807.75
$ws.Range("H132").Value = 1514.75
$ws.Range("I132").Value = 1048.3334
$ws.Range("K132").Value = 9435.000599999999
$ws.Range("M132").Value = -6905.000599999999

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 5900
$ws.Range("J92").Value = 5900
$ws.Range("L92").Value = 5900
$ws.Range("N92").Value = -9644
$ws.Range("H94").Value = 32666.334
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("N94").Value = -21352
$ws.Range("H102").Value = 9064844
$ws.Range("I102").Value = 15875392
$ws.Range("J102").Value = 3105613.5
$ws.Range("K102").Value = 15875392
$ws.Range("L102").Value = 3105613.5
$ws.Range("M102").Value = -15873770
$ws.Range("N102").Value = -3108857.5
$ws.Range("H107").Value = 1026.579
$ws.Range("I107").Value = 957.5714
$ws.Range("J107").Value = 1219.8
$ws.Range("K107").Value = 957.5714
$ws.Range("L107").Value = 1219.8
$ws.Range("M107").Value = 962.4286
$ws.Range("N107").Value = -5059.8
$ws.Range("H122").Value = 525902.1
$ws.Range("I122").Value = 637738.3
$ws.Range("K122").Value = 1913214.9
$ws.Range("M122").Value = -1910764.9
$ws.Range("H132").Value = 3176.6428
$ws.Range("I132").Value = 2853.8157
$ws.Range("K132").Value = 8561.447100000001
$ws.Range("M132").Value = -6031.447100000001

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7241.8887
$ws.Range("I122").Value = 3664.6667
$ws.Range("K122").Value = 10994.0001
$ws.Range("M122").Value = -8544.000100000001

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 30499
$ws.Range("J94").Value = 30499
$ws.Range("L94").Value = 30499
$ws.Range("N94").Value = -32301
$ws.Range("H122").Value = 2434.8823
$ws.Range("I122").Value = 1590
$ws.Range("J122").Value = 3641.8572
$ws.Range("K122").Value = 4770
$ws.Range("L122").Value = 10925.5716
$ws.Range("M122").Value = -2320
$ws.Range("N122").Value = -15825.5716
$ws.Range("H136").Value = 4361.5557
$ws.Range("I136").Value = 2438.9412
$ws.Range("J136").Value = 7630
$ws.Range("K136").Value = 7316.823600000001
$ws.Range("L136").Value = 22890
$ws.Range("M136").Value = -4766.823600000001
$ws.Range("N136").Value = -27990
